$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("169:169").Delete()
